$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# The sheet previously had no real header row (row 1 just duplicated the
# row-2 data). Bring row 1 up to the same header schema used on the other
# sheets, and extend row 2 with the property_category/category/date/
# legislator_name/legislator_id/source_file/index columns.

# Pull in the existing header/data cell formatting for the new columns
# before writing values into them, so F:L match B:E's look (bold+border
# header in row 1, plain in row 2).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("F2:L2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 1 headers
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "property_category"
$ws.Range("G1").Value = "category"
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("K1").Value = "source_file"
$ws.Range("L1").Value = "index"

# Row 2 data - existing columns (name/quantity/owner/total) stay the same,
# fill the new trailing columns.
$ws.Range("F2").Value = "otherbonds"
$ws.Range("G2").Value = "normal"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2013-12-30"
$ws.Range("I2").Value = "蘇震清"
$ws.Range("J2").Value = 1718
$ws.Range("K2").Value = "tmpb3b61"
$ws.Range("L2").Value = 75
